# The workbook has a single worksheet ("Pie Chart") containing a small
# table (A1:B4) that feeds a 3-D pie chart. Update the three data values
# in column B to their new numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 18

$wb.Save()
